$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-08 14:30:46"
$wsZhCn.Range("G2").Value = "2016-01-08 14:31:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-08 14:31:00"
$wsDeDe.Range("G2").Value = "2016-01-08 14:32:16"
